# Applies the "unlock-testing-session.xlsx" edit:
#  - Window/pane/selection view tweaks on the workbook + first sheet
#  - Narrower column J
#  - Updated notes / header text (shared strings content, via cell values)
#  - Row 6/7 "No. Trials/app" (K) corrected from 4 to 1
#  - Row 8 filled in with the 2014-01-21 session data + taller row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unlock Testing ")

# ---- Workbook window settings ---------------------------------------
$wb.Windows.Item(1).Top    = 120
$wb.Windows.Item(1).Height = 9465

# ---- Header text: "Placing cap" -> "Who put cap on?" ----------------
$ws.Range("E4").Value = "Who put cap on?"

# ---- Column J width ----------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 10.28515625

# ---- Row 6 / Row 7: trials/app corrected 4 -> 1 -----------------------
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 1

# ---- Row 8: fill in the 2014-01-21 session -----------------------------
$ws.Range("C8").Value = 20140121
$ws.Range("D8").Value = "20140121-004"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = "-"
$ws.Range("I8").Value = "-"
$ws.Range("J8").Value = "SSVEP"
$ws.Range("K8").Value = "-"
$ws.Range("L8").Value = "-"
$ws.Range("N8").Value = "NA"
$ws.Range("P8").Value = "-"
$ws.Range("Q8").Value = "-"
$ws.Range("R8").Value = "HSD"
$ws.Range("S8").Value = "Planned to used purple cap"
$ws.Range("T8").Value = "James cancelled due to weather conditions expected to happen in the afternoon"
$ws.Range("V8").Value = "-"
$ws.Range("W8").Value = "-"
$ws.Range("X8").Value = "-"

$ws.Rows.Item(8).RowHeight = 45

# ---- Sheet view: scroll position + split pane + selections -----------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 4
$win.SplitRow = 4
$win.FreezePanes = $false
$ws.Range("C8").Select()
$win.Split = $true
$ws.Range("E4").Select()
